{"js": "const body = context.document.body;\nconst pairs = [\n  [\"149\u00d79=\", \"999\u00d73=\"],\n  [\"419\u00d73=\", \"447\u00d78=\"],\n  [\"130\u00d78=\", \"509\u00d78=\"],\n  [\"528\u00d72=\", \"673\u00d76=\"],\n  [\"275\u00d72=\", \"410\u00d77=\"],\n  [\"679\u00d78=\", \"821\u00d76=\"],\n  [\"999\u00d74=\", \"219\u00d75=\"],\n  [\"344\u00d76=\", \"391\u00d74=\"],\n  [\"582\u00d78=\", \"585\u00d72=\"],\n  [\"960\u00d75=\", \"441\u00d79=\"],\n  [\"186\u00d72=\", \"507\u00d75=\"],\n  [\"816\u00d79=\", \"563\u00d76=\"],\n  [\"666\u00d77=\", \"715\u00d78=\"],\n  [\"959\u00d79=\", \"949\u00d76=\"],\n  [\"974\u00d79=\", \"948\u00d76=\"],\n  [\"842\u00d76=\", \"904\u00d72=\"],\n  [\"650\u00d76=\", \"763\u00d76=\"],\n  [\"307\u00d76=\", \"953\u00d77=\"],\n  [\"756\u00d75=\", \"866\u00d72=\"],\n  [\"241\u00d74=\", \"332\u00d79=\"],\n  [\"224\u00d73=\", \"976\u00d77=\"],\n  [\"249\u00d79=\", \"345\u00d74=\"],\n  [\"336\u00d73=\", \"898\u00d74=\"],\n  [\"592\u00d74=\", \"270\u00d75=\"],\n  [\"904\u00d76=\", \"174\u00d76=\"],\n];\n\nfor (const [findText, replaceText] of pairs) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{find=\"149\u00d79=\"; replace=\"999\u00d73=\"}\n    @{find=\"419\u00d73=\"; replace=\"447\u00d78=\"}\n    @{find=\"130\u00d78=\"; replace=\"509\u00d78=\"}\n    @{find=\"528\u00d72=\"; replace=\"673\u00d76=\"}\n    @{find=\"275\u00d72=\"; replace=\"410\u00d77=\"}\n    @{find=\"679\u00d78=\"; replace=\"821\u00d76=\"}\n    @{find=\"999\u00d74=\"; replace=\"219\u00d75=\"}\n    @{find=\"344\u00d76=\"; replace=\"391\u00d74=\"}\n    @{find=\"582\u00d78=\"; replace=\"585\u00d72=\"}\n    @{find=\"960\u00d75=\"; replace=\"441\u00d79=\"}\n    @{find=\"186\u00d72=\"; replace=\"507\u00d75=\"}\n    @{find=\"816\u00d79=\"; replace=\"563\u00d76=\"}\n    @{find=\"666\u00d77=\"; replace=\"715\u00d78=\"}\n    @{find=\"959\u00d79=\"; replace=\"949\u00d76=\"}\n    @{find=\"974\u00d79=\"; replace=\"948\u00d76=\"}\n    @{find=\"842\u00d76=\"; replace=\"904\u00d72=\"}\n    @{find=\"650\u00d76=\"; replace=\"763\u00d76=\"}\n    @{find=\"307\u00d76=\"; replace=\"953\u00d77=\"}\n    @{find=\"756\u00d75=\"; replace=\"866\u00d72=\"}\n    @{find=\"241\u00d74=\"; replace=\"332\u00d79=\"}\n    @{find=\"224\u00d73=\"; replace=\"976\u00d77=\"}\n    @{find=\"249\u00d79=\"; replace=\"345\u00d74=\"}\n    @{find=\"336\u00d73=\"; replace=\"898\u00d74=\"}\n    @{find=\"592\u00d74=\"; replace=\"270\u00d75=\"}\n    @{find=\"904\u00d76=\"; replace=\"174\u00d76=\"}\n)\n\nforeach ($p in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.find\n    $find.Replacement.Text = $p.replace\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
